# Auto-generated edit script: applies cell value updates per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 2016
$ws.Range("I62").Value = 2016
$ws.Range("K62").Value = 2016
$ws.Range("M62").Value = -1392

# Row 65
$ws.Range("H65").Value = 2016
$ws.Range("I65").Value = 2016
$ws.Range("K65").Value = 10080
$ws.Range("M65").Value = -6960

# Row 70
$ws.Range("H70").Value = 3312
$ws.Range("J70").Value = 3666.6667
$ws.Range("L70").Value = 11000.0001
$ws.Range("N70").Value = -11540.0001

# Row 73
$ws.Range("H73").Value = 3312
$ws.Range("J73").Value = 3666.6667
$ws.Range("L73").Value = 11000.0001
$ws.Range("N73").Value = -12872.0001

# Row 76
$ws.Range("H76").Value = 3000
$ws.Range("I76").Value = 3000
$ws.Range("K76").Value = 3000
$ws.Range("M76").Value = -2685

# Row 79
$ws.Range("H79").Value = 3000
$ws.Range("I79").Value = 3000
$ws.Range("K79").Value = 3000
$ws.Range("M79").Value = -1908

# Row 87
$ws.Range("H87").Value = 49997.5
$ws.Range("J87").Value = 49997.5
$ws.Range("L87").Value = 49997.5
$ws.Range("N87").Value = -52493.5

# Row 90
$ws.Range("H90").Value = 49997.5
$ws.Range("J90").Value = 49997.5
$ws.Range("L90").Value = 149992.5
$ws.Range("N90").Value = -162472.5

# Row 92
$ws.Range("H92").Value = 339.47058
$ws.Range("I92").Value = 298.06668
$ws.Range("K92").Value = 298.06668
$ws.Range("M92").Value = 949.93332

# Row 116
$ws.Range("H116").Value = 5500
$ws.Range("I116").Value = 4500
$ws.Range("J116").Value = 7500
$ws.Range("K116").Value = 4500
$ws.Range("L116").Value = 7500
$ws.Range("M116").Value = -1058
$ws.Range("N116").Value = -14384

$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 54785
$ws.Range("J24").Value = 54785
$ws.Range("L24").Value = 54785
$ws.Range("N24").Value = -55533

# Row 32
$ws.Range("H32").Value = 3391.7556
$ws.Range("I32").Value = 1253.119
$ws.Range("K32").Value = 1253.119
$ws.Range("M32").Value = -966.1189999999999

# Row 61
$ws.Range("H61").Value = 2654.5
$ws.Range("J61").Value = 2243.5
$ws.Range("L61").Value = 2243.5
$ws.Range("N61").Value = -2667.5

# Row 80
$ws.Range("H80").Value = 29998.125
$ws.Range("J80").Value = 29998.125
$ws.Range("L80").Value = 29998.125
$ws.Range("N80").Value = -31994.125

# Row 83
$ws.Range("H83").Value = 29998.125
$ws.Range("J83").Value = 29998.125
$ws.Range("L83").Value = 89994.375
$ws.Range("N83").Value = -99978.375

# Row 100
$ws.Range("H100").Value = 54785
$ws.Range("J100").Value = 54785
$ws.Range("L100").Value = 54785
$ws.Range("N100").Value = -56949

# Row 136
$ws.Range("H136").Value = 2654.5
$ws.Range("J136").Value = 2243.5
$ws.Range("L136").Value = 6730.5
$ws.Range("N136").Value = -11830.5

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 8504
$ws.Range("I20").Value = 9008
$ws.Range("K20").Value = 9008
$ws.Range("M20").Value = -8761

# Row 35
$ws.Range("H35").Value = 14999
$ws.Range("J35").Value = 14999
$ws.Range("L35").Value = 14999
$ws.Range("N35").Value = -15619

# Row 82
$ws.Range("H82").Value = 20403.111
$ws.Range("J82").Value = 29714.285
$ws.Range("L82").Value = 29714.285
$ws.Range("N82").Value = -30480.285

# Row 85
$ws.Range("H85").Value = 20403.111
$ws.Range("J85").Value = 29714.285
$ws.Range("L85").Value = 29714.285
$ws.Range("N85").Value = -32366.285

# Row 134
$ws.Range("H134").Value = 5140.4546
$ws.Range("J134").Value = 3663.9375
$ws.Range("L134").Value = 10991.8125
$ws.Range("N134").Value = -16061.8125

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 653.75
$ws.Range("J7").Value = 560
$ws.Range("L7").Value = 560
$ws.Range("N7").Value = -786

# Row 31
$ws.Range("H31").Value = 1303
$ws.Range("I31").Value = 1130.5454
$ws.Range("K31").Value = 1130.5454
$ws.Range("M31").Value = -835.5454

# Row 34
$ws.Range("H34").Value = 1303
$ws.Range("I34").Value = 1130.5454
$ws.Range("K34").Value = 1130.5454
$ws.Range("M34").Value = -928.5454

# Row 41
$ws.Range("H41").Value = 13277.556
$ws.Range("J41").Value = 14562.25
$ws.Range("L41").Value = 14562.25
$ws.Range("N41").Value = -15418.25

# Row 50
$ws.Range("H50").Value = 19536.285
$ws.Range("J50").Value = 18926.428
$ws.Range("L50").Value = 18926.428
$ws.Range("N50").Value = -20176.428

# Row 51
$ws.Range("H51").Value = 20000
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

# Row 59
$ws.Range("H59").Value = 27450.875
$ws.Range("I59").Value = 19804
$ws.Range("J59").Value = 29999.834
$ws.Range("K59").Value = 19804
$ws.Range("L59").Value = 29999.834
$ws.Range("M59").Value = -18659
$ws.Range("N59").Value = -32289.834

# Row 60
$ws.Range("H60").Value = 19631.777
$ws.Range("I60").Value = 24022.75
$ws.Range("J60").Value = 16119
$ws.Range("K60").Value = 24022.75
$ws.Range("L60").Value = 16119
$ws.Range("M60").Value = -23511.75
$ws.Range("N60").Value = -17141

# Row 61
$ws.Range("H61").Value = 20000
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

# Row 68
$ws.Range("H68").Value = 28180.363
$ws.Range("J68").Value = 29998.334
$ws.Range("L68").Value = 29998.334
$ws.Range("N68").Value = -31496.334

# Row 71
$ws.Range("H71").Value = 28180.363
$ws.Range("J71").Value = 29998.334
$ws.Range("L71").Value = 89995.00199999999
$ws.Range("N71").Value = -97483.00199999999

# Row 74
$ws.Range("H74").Value = 29998.4
$ws.Range("J74").Value = 29998.4
$ws.Range("L74").Value = 29998.4
$ws.Range("N74").Value = -31746.4

# Row 77
$ws.Range("H77").Value = 29998.4
$ws.Range("J77").Value = 29998.4
$ws.Range("L77").Value = 89995.20000000001
$ws.Range("N77").Value = -98731.20000000001

# Row 96
$ws.Range("H96").Value = 23762
$ws.Range("J96").Value = 23762
$ws.Range("L96").Value = 23762
$ws.Range("N96").Value = -29254

$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 562.4
$ws.Range("I122").Value = 478
$ws.Range("J122").Value = 689
$ws.Range("K122").Value = 4302
$ws.Range("L122").Value = 6201
$ws.Range("M122").Value = -1852
$ws.Range("N122").Value = -11101

$ws = $wb.Worksheets.Item("LTW")
# Row 69
$ws.Range("H69").Value = 70000
$ws.Range("J69").Value = 70000
$ws.Range("L69").Value = 70000
$ws.Range("N69").Value = -71622

# Row 72
$ws.Range("H72").Value = 70000
$ws.Range("J72").Value = 70000
$ws.Range("L72").Value = 210000
$ws.Range("N72").Value = -218112

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 5511.8335
$ws.Range("I136").Value = 3394.5
$ws.Range("J136").Value = 8158.5
$ws.Range("K136").Value = 10183.5
$ws.Range("L136").Value = 24475.5
$ws.Range("M136").Value = -7633.5
$ws.Range("N136").Value = -29575.5
